# Apply updated cryptocurrency price/volume data to the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.486.21'
$ws.Range("E2").Value = '  -0.05%  '
$ws.Range("D3").Value = '1.567.90'
$ws.Range("E3").Value = '  -2.11%  '
$ws.Range("D5").Value = "'211.90"
$ws.Range("E5").Value = '  -1.45%  '
$ws.Range("E6").Value = '  -1.15%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").Value = "'46.03"
$ws.Range("E8").Value = '  +4.36%  '
$ws.Range("D9").Value = "'24.02"
$ws.Range("E9").Value = '  +0.01%  '
$ws.Range("E10").Value = '  -1.87%  '
$ws.Range("E11").Value = '  -1.56%  '
$ws.Range("E12").Value = '  -0.32%  '
$ws.Range("D13").Value = '1.791.43'
$ws.Range("E13").Value = '  -2.15%  '
$ws.Range("D14").Value = '1.562.79'
$ws.Range("E14").Value = '  -2.47%  '
$ws.Range("E15").Value = '  -2.67%  '
$ws.Range("E16").Value = '  -2.98%  '
$ws.Range("D17").Value = '28.499.72'
$ws.Range("E17").Value = '  +0.00%  '
$ws.Range("D18").Value = "'62.25"
$ws.Range("E18").Value = '  -1.78%  '
$ws.Range("D19").Value = "'227.65"
$ws.Range("E19").Value = '  -1.86%  '
$ws.Range("D20").Value = "'7.35"
$ws.Range("E20").Value = '  -2.60%  '
$ws.Range("E21").Value = '  -2.75%  '
$ws.Range("E23").Value = '  -5.92%  '
$ws.Range("D24").Value = "'9.14"
$ws.Range("E24").Value = '  -3.01%  '
$ws.Range("E25").Value = '  +6.20%  '
$ws.Range("D26").Value = "'151.08"
$ws.Range("E26").Value = '  -0.95%  '
$ws.Range("D27").Value = "'15.01"
$ws.Range("E27").Value = '  -2.03%  '
$ws.Range("E28").Value = '  -2.72%  '
$ws.Range("E29").Value = '  -4.09%  '
$ws.Range("E31").Value = '  -1.75%  '
$ws.Range("E32").Value = '  -4.21%  '
$ws.Range("E33").Value = '  -1.38%  '
$ws.Range("D34").Value = "'3.07"
$ws.Range("E34").Value = '  -2.95%  '
$ws.Range("D35").Value = '1.396.37'
$ws.Range("E35").Value = '  -1.93%  '
$ws.Range("E36").Value = '  -0.79%  '
$ws.Range("E37").Value = '  -3.46%  '
$ws.Range("E38").Value = '  +1.22%  '
$ws.Range("D39").Value = "'2.59"
$ws.Range("E39").Value = '  +2.62%  '
$ws.Range("E40").Value = '  -1.06%  '
$ws.Range("D41").Value = "'0.536"
$ws.Range("E41").Value = '  -1.49%  '
$ws.Range("B43").Value = 'ARBITRUM'
$ws.Range("C43").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D43").Value = "'0.786"
$ws.Range("E43").Value = '  -4.52%  '
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").Value = "'1.88"
$ws.Range("E44").Value = '  +1.81%  '
$ws.Range("E45").Value = '  -4.27%  '
$ws.Range("E46").Value = '  -1.09%  '
$ws.Range("D47").Value = "'62.91"
$ws.Range("E47").Value = '  -3.16%  '
$ws.Range("D48").Value = '1.704.45'
$ws.Range("E48").Value = '  -2.12%  '
$ws.Range("D49").Value = "'85.90"
$ws.Range("E49").Value = '  -1.87%  '
$ws.Range("D50").Value = "'0.0526"
$ws.Range("E50").Value = '  -0.27%  '
$ws.Range("E51").Value = '  -4.70%  '
